$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "77÷6=12, 5"
$t.Cell(1,2).Range.Text = "50÷7=7, 1"
$t.Cell(1,3).Range.Text = "11÷2=5, 1"
$t.Cell(1,4).Range.Text = "94÷2=47, 0"
$t.Cell(1,5).Range.Text = "54÷4=13, 2"

$t.Cell(5,1).Range.Text = "49÷6=8, 1"
$t.Cell(5,2).Range.Text = "14÷9=1, 5"
$t.Cell(5,3).Range.Text = "42÷8=5, 2"
$t.Cell(5,4).Range.Text = "41÷5=8, 1"
$t.Cell(5,5).Range.Text = "99÷5=19, 4"

$t.Cell(9,1).Range.Text = "93÷5=18, 3"
$t.Cell(9,2).Range.Text = "94÷8=11, 6"
$t.Cell(9,3).Range.Text = "40÷5=8, 0"
$t.Cell(9,4).Range.Text = "79÷4=19, 3"
$t.Cell(9,5).Range.Text = "13÷3=4, 1"

$t.Cell(13,1).Range.Text = "96÷2=48, 0"
$t.Cell(13,2).Range.Text = "33÷4=8, 1"
$t.Cell(13,3).Range.Text = "13÷3=4, 1"
$t.Cell(13,4).Range.Text = "17÷7=2, 3"
$t.Cell(13,5).Range.Text = "75÷9=8, 3"

$t.Cell(17,1).Range.Text = "64÷8=8, 0"
$t.Cell(17,2).Range.Text = "21÷5=4, 1"
$t.Cell(17,3).Range.Text = "71÷7=10, 1"
$t.Cell(17,4).Range.Text = "32÷2=16, 0"
$t.Cell(17,5).Range.Text = "51÷9=5, 6"

